$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.214.13"
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = "'1.570.91"
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = "'211.21"
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").Value = "'22.07"
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D11").Value = "'0.0870"
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = "'1.563.79"
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = "'3.79"
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = "'27.156.27"
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = "'62.24"
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = "'7.49"
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = "'0.0₃0703"
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'216.41"
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").Value = "'153.96"
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = "'15.08"
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  +1.89%  '
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("D32").Value = "'3.25"
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("D33").Value = "'1.450.29"
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("E34").Value = '  +1.51%  '
$ws.Range("E35").Value = '  +7.43%  '
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").Value = "'2.35"
$ws.Range("E37").Value = '  +0.73%  '
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").Value = "'5.86"
$ws.Range("E40").Value = '  +2.45%  '
$ws.Range("D41").Value = "'0.811"
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("D45").Value = "'64.45"
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("D47").Value = "'1.705.39"
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").Value = "'86.17"
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("E49").Value = '  +3.66%  '
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = "'0.0958"
$ws.Range("E51").Value = '  +0.24%  '
